# Paste a new weld-record row into the template (row 15), pushing the
# template's trailing blank spacer row down to row 16 - mirroring how the
# existing rows 2-14 were built (boilerplate pipe/spec/material/method text
# reused from the row above, a fresh seam number + wall thickness + preheat
# typed in, and the weld date copied from an existing dated row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at 15; this shifts the old (empty) row 15 down to
# row 16, preserving it as the sheet's trailing spacer row.
$ws.Rows.Item(15).Insert()

# Templated (boilerplate) columns: reuse row 14's content + exact formatting.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("H15").PasteSpecial(-4122)

# Date column: copy an earlier row's formatting (matches rows 3-13, not 14).
$ws.Range("J13").Copy()
$ws.Range("J15").PasteSpecial(-4122)

# Freshly typed data columns (weld no., spec, preheat): plain worksheet
# default formatting, picked up from an untouched cell in the same columns.
$ws.Range("B200").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C200").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("I200").Copy()
$ws.Range("I15").PasteSpecial(-4122)

$ws.Range("A15").Value = "2310-P-02161-B62SA0"
$ws.Range("B15").Value = "11"
$ws.Range("C15").Value = "4270"
$ws.Range("D15").Value = "323.9*6.35"
$ws.Range("E15").Value = "A312GR.TP304"
$ws.Range("F15").Value = "1G"
$ws.Range("G15").Value = "GTAW+SMAW"
$ws.Range("H15").Value = "ER308/E308-16"
$ws.Range("I15").Value = "16"
$ws.Range("J15").Value = "2021/10/27"

$ws.Rows.Item(15).RowHeight = 27

$ws.Range("D12").Select()
